$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing D:K data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy the number/date formatting from column F (the old column D, now shifted)
# into the two new columns D and E so the new cells render identically to their neighbours.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("F7:F102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new columns with the latest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 333300
$ws.Range("E8").Value = 381500
$ws.Range("D9").Value = 252900
$ws.Range("E9").Value = 283600
$ws.Range("D10").Value = 80400
$ws.Range("E10").Value = 97900
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 258800
$ws.Range("E17").Value = 289700
$ws.Range("D18").Value = 74500
$ws.Range("E18").Value = 91800
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 104300
$ws.Range("E21").Value = 124100
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 74500
$ws.Range("E23").Value = 91800
$ws.Range("D24").Value = 16800
$ws.Range("E24").Value = 19200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 57700
$ws.Range("E26").Value = 72600
$ws.Range("D27").Value = 57700
$ws.Range("E27").Value = 72600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 57700
$ws.Range("E33").Value = 72600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 57700
$ws.Range("E35").Value = 72600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 17100
$ws.Range("E41").Value = 10000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 134200
$ws.Range("E43").Value = 180500
$ws.Range("D44").Value = 251300
$ws.Range("E44").Value = 238900
$ws.Range("D45").Value = 7000
$ws.Range("E45").Value = 7800
$ws.Range("D46").Value = 409500
$ws.Range("E46").Value = 437100
$ws.Range("D47").Value = 65000
$ws.Range("E47").Value = 63600
$ws.Range("D48").Value = 1627200
$ws.Range("E48").Value = 1624700
$ws.Range("D49").Value = 236900
$ws.Range("E49").Value = 237700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 16800
$ws.Range("E52").Value = 16300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2355400
$ws.Range("E54").Value = 2379500
$ws.Range("D57").Value = 77600
$ws.Range("E57").Value = 92500
$ws.Range("D58").Value = 36500
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 66900
$ws.Range("E59").Value = 62200
$ws.Range("D60").Value = 181000
$ws.Range("E60").Value = 154700
$ws.Range("D61").Value = 589900
$ws.Range("E61").Value = 631300
$ws.Range("D62").Value = 164100
$ws.Range("E62").Value = 161000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 935100
$ws.Range("E66").Value = 946900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1423700
$ws.Range("E72").Value = 1428300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1420300
$ws.Range("E76").Value = 1432600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 57700
$ws.Range("E81").Value = 72600
$ws.Range("D83").Value = 29800
$ws.Range("E83").Value = 32300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 119200
$ws.Range("E89").Value = 93000
$ws.Range("D91").Value = -33000
$ws.Range("E91").Value = -40400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -33000
$ws.Range("E94").Value = -40400
$ws.Range("D96").Value = -4700
$ws.Range("E96").Value = -4800
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -79200
$ws.Range("E100").Value = -95700
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 7100
$ws.Range("E102").Value = -43100
